$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122 (-4122 == 10 as unsigned; used with PasteSpecial)
$xlPasteFormats     = -4122
$xlRight            = -4152
$xlBottom           = -4107

# --- Column P header (row 4): reuse the exact same style as O4 (plain year, s=1) ---
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial($xlPasteFormats)
$ws.Range("P4").Value = 2022

# --- P5 (totals row): same look as O5 but with new number format "#,##0.0" and
#     right alignment only (no vertical centering). Build the target style once on
#     a scratch cell so only a single new cellXfs entry is produced. ---
$ws.Range("O5").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$style = $ws.Range("Z1").Style
$style.NumberFormat = "#,##0.0"
$style.HorizontalAlignment = $xlRight
$style.VerticalAlignment = $xlBottom
$ws.Range("Z1").Copy()
$ws.Range("P5").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Clear()
$ws.Range("P5").Value = 1188.7

# --- P6:P15 (regular data rows): same treatment, built once on a scratch cell and
#     pasted onto the whole block so only one new cellXfs entry is produced. ---
$ws.Range("O6").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$style2 = $ws.Range("Z1").Style
$style2.NumberFormat = "#,##0.0"
$style2.HorizontalAlignment = $xlRight
$style2.VerticalAlignment = $xlBottom
$ws.Range("Z1").Copy()
$ws.Range("P6:P15").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Clear()

$ws.Range("P6").Value = 263.89999999999998
$ws.Range("P7").Value = 263.2
$ws.Range("P8").Value = 12.4
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = 93
$ws.Range("P11").Value = 171.5
$ws.Range("P12").Value = 220.6
$ws.Range("P13").Value = 159.30000000000001
$ws.Range("P14").Value = 1.7
$ws.Range("P15").Value = "-"

# --- P16 (bottom totals row, with border): same treatment again via scratch cell. ---
$ws.Range("O16").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$style3 = $ws.Range("Z1").Style
$style3.NumberFormat = "#,##0.0"
$style3.HorizontalAlignment = $xlRight
$style3.VerticalAlignment = $xlBottom
$ws.Range("Z1").Copy()
$ws.Range("P16").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Clear()
$ws.Range("P16").Value = 3.1

# --- Match the recorded cursor position in the sheet view ---
$ws.Range("Q7").Select()
